$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.321.29'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '1.880.79'
$ws.Range('E3').Value = '  -1.04%  '
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').Value = "'245.39"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.83%  '
$ws.Range('D6').Value = "'0.679"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.15%  '
$ws.Range('E7').Value = '  -0.64%  '
$ws.Range('E8').Value = '  +5.09%  '
$ws.Range('D9').Value = "'0.358"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('D10').Value = "'53.24"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('D11').Value = "'0.0741"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.15%  '
$ws.Range('D12').Value = "'0.0978"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = "'13.51"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.59%  '
$ws.Range('D14').Value = '2.150.57'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = "'0.768"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.95%  '
$ws.Range('D16').Value = "'4.93"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').Value = '1.820.62'
$ws.Range('E17').Value = '  -4.18%  '
$ws.Range('D18').Value = '35.293.49'
$ws.Range('E18').Value = '  +0.44%  '
$ws.Range('D19').Value = "'73.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('D20').Value = '0.0₃0823'
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D21').Value = "'244.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.06%  '
$ws.Range('D22').Value = "'12.86"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.31%  '
$ws.Range('D23').Value = "'5.05"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Value = "'2.70"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.15%  '
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('D26').Value = "'2.16"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.48%  '
$ws.Range('D27').Value = "'164.95"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.14%  '
$ws.Range('D28').Value = "'8.61"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.28%  '
$ws.Range('D29').Value = "'18.28"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.21%  '
$ws.Range('E30').Value = '  -1.99%  '
$ws.Range('D31').Value = "'4.29"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('D32').Value = "'0.0592"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.67%  '
$ws.Range('D33').Value = "'4.18"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.40%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = "'1.85"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -11.11%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').Value = "'1.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').Value = "'1.42"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -13.44%  '
$ws.Range('D37').Value = "'0.852"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.21%  '
$ws.Range('D38').Value = "'1.95"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.05%  '
$ws.Range('D39').Value = "'0.0730"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +8.84%  '
$ws.Range('D40').Value = "'17.36"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.19%  '
$ws.Range('D41').Value = "'0.0218"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('D42').Value = "'96.76"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.08%  '
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('D44').Value = '1.308.05'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('D45').Value = "'2.39"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.14%  '
$ws.Range('D46').Value = "'0.0799"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.92%  '
$ws.Range('E47').Value = '  -1.17%  '
$ws.Range('E48').Value = '  -0.82%  '
$ws.Range('D49').Value = "'11.81"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.76%  '
$ws.Range('D50').Value = "'6.29"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.87%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = "'42.02"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.61%  '
